# ShearF-HW50 : "Updated notebook, reran simulation"
#
# The simulation notebook added two brand-new HKL methods ("Holden" and
# "Rizzie Spiral"). They get inserted into the method list right after
# "Spiral5", which is why the table now has two extra rows (4 and 5)
# holding their freshly simulated numbers, while every method that used
# to start at row 4 ("RotRing OmegaMax-90" and onward) is pushed down by
# two rows, carrying its previously simulated numbers with it. Two more
# rows (30, 31) are appended at the bottom so the table keeps covering
# every method in the (now longer) list. One existing label was also
# corrected along the way: "Thomas Hex" -> "Matthies Hex".
#
# Column A/B simply keep counting 0,1,2,... down the method list (same
# rule as before, just two rows longer) - only the simulated payload in
# columns C:T needs to move/appear/get-appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Slide the existing simulation payload (columns C:T only) down by
#    two rows: old row 29 -> 31, 28 -> 30, ..., 4 -> 6. Walk from the
#    bottom up so a source row is never overwritten before it's been
#    copied from. Row.Copy also means no stray new style entries get
#    created for the numeric cells.
# ------------------------------------------------------------------
for ($r = 29; $r -ge 4; $r--) {
    $dest = $r + 2
    $src = $ws.Range("C" + $r + ":T" + $r)
    $trg = $ws.Range("C" + $dest + ":T" + $dest)
    $src.Copy($trg)
}

# ------------------------------------------------------------------
# 2) Fill in the freshly simulated numbers for the two new methods,
#    "Holden" (row 4) and "Rizzie Spiral" (row 5).
# ------------------------------------------------------------------
$ws.Range("C4").Value2 = 1.002336354629532
$ws.Range("D4").Value2 = 0.9994159089285647
$ws.Range("E4").Value2 = 1.002336354629532
$ws.Range("F4").Value2 = 0.9994159089285647
$ws.Range("G4").Value2 = 0.9994159089285647
$ws.Range("H4").Value2 = 1.001284993553495
$ws.Range("I4").Value2 = 0.9984424292754448
$ws.Range("J4").Value2 = 0.9994159089285647
$ws.Range("K4").Value2 = 0.9994159089285647
$ws.Range("L4").Value2 = 0.9994159089285647
$ws.Range("M4").Value2 = 1.000876131779048
$ws.Range("N4").Value2 = 1.000876131779048
$ws.Range("O4").Value2 = 1.001012419037197
$ws.Range("P4").Value2 = 1.000389390828887
$ws.Range("Q4").Value2 = 1.000389390828887
$ws.Range("R4").Value2 = 1.000146020353806
$ws.Range("S4").Value2 = 1.000146020353806
$ws.Range("T4").Value2 = 1.000051917374028

$ws.Range("C5").Value2 = 0.9975845963308388
$ws.Range("D5").Value2 = 1.000603841035816
$ws.Range("E5").Value2 = 0.9975845963308388
$ws.Range("F5").Value2 = 1.000603841035816
$ws.Range("G5").Value2 = 1.000603841035816
$ws.Range("H5").Value2 = 0.9986715319042414
$ws.Range("I5").Value2 = 1.00161027037485
$ws.Range("J5").Value2 = 1.000603841035816
$ws.Range("K5").Value2 = 1.000603841035816
$ws.Range("L5").Value2 = 1.000603841035816
$ws.Range("M5").Value2 = 0.9990942186833276
$ws.Range("N5").Value2 = 0.9990942186833276
$ws.Range("O5").Value2 = 0.9989533230902988
$ws.Range("P5").Value2 = 0.9995974261341573
$ws.Range("Q5").Value2 = 0.9995974261341573
$ws.Range("R5").Value2 = 0.9998490298595719
$ws.Range("S5").Value2 = 0.9998490298595719
$ws.Range("T5").Value2 = 0.9999463202862299

# ------------------------------------------------------------------
# 3) Append the two new rows' C:T payload at the bottom (these are the
#    methods that used to be the table's last two rows, now pushed to
#    rows 30/31 by the same "+2" shift as step 1, just spelled out
#    explicitly since the loop above stopped at row 29/31).
# ------------------------------------------------------------------
# (rows 30 and 31 were already populated by the loop above: old row 28
#  -> new row 30, old row 29 -> new row 31)

# ------------------------------------------------------------------
# 4) Column A/B: extend the sequential 0,1,2,... index down through the
#    two new rows, and make sure every label reads correctly against
#    the longer method list (two new entries inserted after "Spiral5",
#    "Thomas Hex" renamed to "Matthies Hex"). Style is copied from the
#    row above so column A keeps its existing bordered/bold format with
#    no new style table entries.
# ------------------------------------------------------------------
$ws.Range("A29:B29").Copy($ws.Range("A30:B30"))
$ws.Range("A29:B29").Copy($ws.Range("A31:B31"))

$ws.Range("A4").Value2  = 2
$ws.Range("B4").Value2  = "Holden"
$ws.Range("A5").Value2  = 3
$ws.Range("B5").Value2  = "Rizzie Spiral"
$ws.Range("A6").Value2  = 4
$ws.Range("B6").Value2  = "RotRing OmegaMax-90"
$ws.Range("A7").Value2  = 5
$ws.Range("B7").Value2  = "Equal Angle"
$ws.Range("A8").Value2  = 6
$ws.Range("B8").Value2  = "Tilt Rotate"
$ws.Range("A9").Value2  = 7
$ws.Range("B9").Value2  = "CLR"
$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = "Rizzie Hex"
$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = "Matthies Hex"
$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = "Tilt Rotate_Partial"
$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "RotRing OmegaMax-60"
$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = "Equal Angle_Partial"
$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = "Rizzie Hex_Partial"
$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "ND Single"
$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "RD Single"
$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "TD Single"
$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "Morris Single"
$ws.Range("A20").Value2 = 18
$ws.Range("B20").Value2 = "Ring Perpendicular to ND"
$ws.Range("A21").Value2 = 19
$ws.Range("B21").Value2 = "Ring Perpendicular to RD"
$ws.Range("A22").Value2 = 20
$ws.Range("B22").Value2 = "Ring Perpendicular to TD"
$ws.Range("A23").Value2 = 21
$ws.Range("B23").Value2 = "OffsetFTD"
$ws.Range("A24").Value2 = 22
$ws.Range("B24").Value2 = "OffsetATD"
$ws.Range("A25").Value2 = 23
$ws.Range("B25").Value2 = "OffsetF45"
$ws.Range("A26").Value2 = 24
$ws.Range("B26").Value2 = "OffsetA45"
$ws.Range("A27").Value2 = 25
$ws.Range("B27").Value2 = "OffsetFRD"
$ws.Range("A28").Value2 = 26
$ws.Range("B28").Value2 = "OffsetARD"
$ws.Range("A29").Value2 = 27
$ws.Range("B29").Value2 = "Gaussian Quadrature"
$ws.Range("A30").Value2 = 28
$ws.Range("B30").Value2 = "Michael-CCHex"
$ws.Range("A31").Value2 = 29
$ws.Range("B31").Value2 = "Michael-SNHex"
